$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 5
$ws.Range("G6").Value = 4
$ws.Range("G7").Value = 4
$ws.Range("G8").Value = 2
$ws.Range("G9").Value = 5
$ws.Range("G10").Value = 2
$ws.Range("G11").Value = 4
$ws.Range("G12").Value = 3
$ws.Range("G13").Value = 3
$ws.Range("G14").Value = 5
$ws.Range("G15").Value = 1
$ws.Range("G16").Value = 5
$ws.Range("G17").Value = 9
$ws.Range("G18").Value = 5
$ws.Range("G19").Value = 1
$ws.Range("G20").Value = 6
$ws.Range("G21").Value = 4
$ws.Range("G22").Value = 6
$ws.Range("G23").Value = 4
$ws.Range("G24").Value = 1
$ws.Range("G25").Value = 1
$ws.Range("G26").Value = 2
$ws.Range("G27").Value = 2
$ws.Range("G28").Value = 2
$ws.Range("G29").Value = 4
$ws.Range("G30").Value = 2
$ws.Range("G31").Value = 5
$ws.Range("G32").Value = 3
$ws.Range("G33").Value = 1
$ws.Range("G34").Value = 2
$ws.Range("G35").Value = 4
$ws.Range("G36").Value = 6
$ws.Range("G37").Value = 2
$ws.Range("G38").Value = 1
